$d = $word.ActiveDocument

# The document currently ends with:
#   1: "IPCV Coursework - Part 2 - Report"
#   2: "Question 3"
#   3: "cv2.HoughCircles is dependent on ..."
#   4: "" (the paragraph with w:spacing before=0 after=159 - the last, empty paragraph)
#
# We need to insert two new paragraphs right after paragraph 3 and before the
# trailing (previously empty) paragraph, and also add text to that trailing
# paragraph.

$anchor = $d.Paragraphs(3)

# Create two fresh, empty paragraphs right after paragraph 3. InsertParagraphAfter
# copies paragraph 3's paragraph formatting (no w:spacing override, jc=start,
# b/bCs/u all "off"), which is exactly what both new paragraphs need.
$null = $anchor.Range.InsertParagraphAfter()
$anchor2 = $d.Paragraphs(3)
$null = $anchor2.Range.InsertParagraphAfter()

# Paragraph 4: bold + underlined "Introduction" heading line.
$introXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="Normal"/>
<w:bidi w:val="0"/>
<w:jc w:val="start"/>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:u w:val="single"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>
<w:b/>
<w:bCs/>
<w:u w:val="single"/>
</w:rPr>
<w:t>Introduction- Provide a brief overview of the assignment, the tasks involved and what you discovered.</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$p4 = $d.Paragraphs(4)
$null = $p4.Range.InsertXML($introXml)

# Paragraph 5: normal-weight body paragraph - plain text assignment is enough
# since the freshly-inserted empty paragraph already carries the matching
# (non-bold, non-underlined) run/paragraph-mark formatting.
$p5 = $d.Paragraphs(5)
$p5.Range.Text = "The assignment involved reconstructing 3D spheres from two images of the scene taken by two virtual cameras with known relative pose (i.e. the setup is calibrated). The solution to the assignment obtains from just the two images, and the known position of the cameras, the 3D centre positions and radius lengths of the spheres, which can then be used to draw spheres to a 3D scene in order to visualise the scene."

# Paragraph 6 is the original trailing (empty) paragraph; it already has the
# correct paragraph/run formatting, so only its text needs to be filled in.
$p6 = $d.Paragraphs(6)
$p6.Range.Text = "Completing the assignment required the use of the epipolar line constraint equation, which was used to compute epipolar lines and match corresponding points. It also involved using the formula for 3D reconstruction to obtain a 3D point from two corresponding 2D image points."
